$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1528053333333333
$ws.Range("N2").Value = 0.458416
$ws.Range("O2").Value = 0.01103433215988526
$ws.Range("P2").Value = 0.01103433215988526
$ws.Range("Q2").Value = 1.440527966453333
$ws.Range("R2").Value = 12.96475169808
$ws.Range("S2").Value = 0.002414404407789243
$ws.Range("T2").Value = 0.002414404407789243

$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("O3").Value = 0.8539197603380489
$ws.Range("P3").Value = 0.8539197603380488
$ws.Range("Q3").Value = 111.47890765388
$ws.Range("R3").Value = 1003.31016888492
$ws.Range("S3").Value = 0.1868448043238856
$ws.Range("T3").Value = 0.1868448043238855

$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 1.712948333333333
$ws.Range("N4").Value = 5.138845
$ws.Range("O4").Value = 0.1236949029880405
$ws.Range("P4").Value = 0.1236949029880405
$ws.Range("Q4").Value = 16.14832365748333
$ws.Range("R4").Value = 145.33491291735
$ws.Range("S4").Value = 0.02706548204893745
$ws.Range("T4").Value = 0.02706548204893745

$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1571906666666667
$ws.Range("N5").Value = 0.471572
$ws.Range("O5").Value = 0.01135100451402528
$ws.Range("P5").Value = 0.01135100451402528
$ws.Range("Q5").Value = 1.481869424706667
$ws.Range("R5").Value = 13.33682482236
$ws.Range("S5").Value = 0.002483694974411864
$ws.Range("T5").Value = 0.002483694974411863

$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1528053333333333
$ws.Range("N6").Value = 0.458416
$ws.Range("O6").Value = 0.01103433215988526
$ws.Range("P6").Value = 0.01103433215988526
$ws.Range("Q6").Value = 2.507586098565333
$ws.Range("R6").Value = 22.568274887088
$ws.Range("S6").Value = 0.004202852752795413
$ws.Range("T6").Value = 0.004202852752795413

$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("O7").Value = 0.8539197603380489
$ws.Range("P7").Value = 0.8539197603380488
$ws.Range("S7").Value = 0.3252484122646248
$ws.Range("T7").Value = 0.3252484122646247

$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 1.712948333333333
$ws.Range("N8").Value = 5.138845
$ws.Range("O8").Value = 0.1236949029880405
$ws.Range("P8").Value = 0.1236949029880405
$ws.Range("Q8").Value = 28.11004913589833
$ws.Range("R8").Value = 252.990442223085
$ws.Range("S8").Value = 0.04711399439469597
$ws.Range("T8").Value = 0.04711399439469596

$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1571906666666667
$ws.Range("N9").Value = 0.471572
$ws.Range("O9").Value = 0.01135100451402528
$ws.Range("P9").Value = 0.01135100451402528
$ws.Range("Q9").Value = 2.579550870110667
$ws.Range("R9").Value = 23.215957830996
$ws.Range("S9").Value = 0.00432346968330346
$ws.Range("T9").Value = 0.004323469683303459

$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1528053333333333
$ws.Range("N10").Value = 0.458416
$ws.Range("O10").Value = 0.01103433215988526
$ws.Range("P10").Value = 0.01103433215988526
$ws.Range("Q10").Value = 1.102194190458667
$ws.Range("R10").Value = 9.919747714127999
$ws.Range("S10").Value = 0.001847338318765858
$ws.Range("T10").Value = 0.001847338318765858

$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("O11").Value = 0.8539197603380489
$ws.Range("P11").Value = 0.8539197603380488
$ws.Range("Q11").Value = 85.29609090290801
$ws.Range("R11").Value = 767.664818126172
$ws.Range("S11").Value = 0.1429609578147989
$ws.Range("T11").Value = 0.1429609578147989

$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 1.712948333333333
$ws.Range("N12").Value = 5.138845
$ws.Range("O12").Value = 0.1236949029880405
$ws.Range("P12").Value = 0.1236949029880405
$ws.Range("Q12").Value = 12.35560081818167
$ws.Range("R12").Value = 111.200407363635
$ws.Range("S12").Value = 0.02070866916228564
$ws.Range("T12").Value = 0.02070866916228564

$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1571906666666667
$ws.Range("N13").Value = 0.471572
$ws.Range("O13").Value = 0.01135100451402528
$ws.Range("P13").Value = 0.01135100451402528
$ws.Range("Q13").Value = 1.133825867297333
$ws.Range("R13").Value = 10.204432805676
$ws.Range("S13").Value = 0.001900354755630373
$ws.Range("T13").Value = 0.001900354755630373

$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1528053333333333
$ws.Range("N14").Value = 0.458416
$ws.Range("O14").Value = 0.01103433215988526
$ws.Range("P14").Value = 0.01103433215988526
$ws.Range("Q14").Value = 1.533205267016889
$ws.Range("R14").Value = 13.798847403152
$ws.Range("S14").Value = 0.002569736680534748
$ws.Range("T14").Value = 0.002569736680534748

$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("O15").Value = 0.8539197603380489
$ws.Range("P15").Value = 0.8539197603380488
$ws.Range("Q15").Value = 118.6509754455054
$ws.Range("R15").Value = 1067.858779009548
$ws.Range("S15").Value = 0.1988655859347398
$ws.Range("T15").Value = 0.1988655859347398

$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 1.712948333333333
$ws.Range("N16").Value = 5.138845
$ws.Range("O16").Value = 0.1236949029880405
$ws.Range("P16").Value = 0.1236949029880405
$ws.Range("Q16").Value = 17.18723652835723
$ws.Range("R16").Value = 154.685128755215
$ws.Range("S16").Value = 0.02880675738212145
$ws.Range("T16").Value = 0.02880675738212145

$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1571906666666667
$ws.Range("N17").Value = 0.471572
$ws.Range("O17").Value = 0.01135100451402528
$ws.Range("P17").Value = 0.01135100451402528
$ws.Range("Q17").Value = 1.577206454787111
$ws.Range("R17").Value = 14.194858093084
$ws.Range("S17").Value = 0.002643485100679584
$ws.Range("T17").Value = 0.002643485100679584
